$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 393.78262
$ws.Range("I41").Value = 151.625
$ws.Range("J41").Value = 522.93335
$ws.Range("K41").Value = 151.625
$ws.Range("L41").Value = 522.93335
$ws.Range("M41").Value = 288.375
$ws.Range("N41").Value = -1402.93335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5177.963
$ws.Range("I86").Value = 6650.25
$ws.Range("J86").Value = 4921.913
$ws.Range("K86").Value = 6650.25
$ws.Range("L86").Value = 4921.913
$ws.Range("M86").Value = -5527.25
$ws.Range("N86").Value = -7167.913

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5177.963
$ws.Range("I89").Value = 6650.25
$ws.Range("J89").Value = 4921.913
$ws.Range("K89").Value = 33251.25
$ws.Range("L89").Value = 24609.565
$ws.Range("M89").Value = -27635.25
$ws.Range("N89").Value = -35841.565

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1642.5883
$ws.Range("I100").Value = 1432.75
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1432.75
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -891.75
$ws.Range("N100").Value = -6082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 2875.25
$ws.Range("I30").Value = 500.33334
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 500.33334
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = -350.33334
$ws.Range("N30").Value = -10300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1464.8125
$ws.Range("I45").Value = 1538.3572
$ws.Range("J45").Value = 950
$ws.Range("K45").Value = 1538.3572
$ws.Range("L45").Value = 950
$ws.Range("M45").Value = -1161.3572
$ws.Range("N45").Value = -1704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10370.182
$ws.Range("I61").Value = 10370.182
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10370.182
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -10158.182
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1013.75
$ws.Range("I74").Value = 808.381
$ws.Range("J74").Value = 1629.8572
$ws.Range("K74").Value = 808.381
$ws.Range("L74").Value = 1629.8572
$ws.Range("M74").Value = 65.61900000000003
$ws.Range("N74").Value = -3377.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1013.75
$ws.Range("I77").Value = 808.381
$ws.Range("J77").Value = 1629.8572
$ws.Range("K77").Value = 4041.905
$ws.Range("L77").Value = 8149.286
$ws.Range("M77").Value = 326.0950000000003
$ws.Range("N77").Value = -16885.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 44126.07
$ws.Range("J121").Value = 44126.07
$ws.Range("L121").Value = 44126.07
$ws.Range("N121").Value = -47620.07

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10370.182
$ws.Range("I136").Value = 10370.182
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 31110.546
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -28560.546
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2793.2
$ws.Range("I86").Value = 2793.2
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2793.2
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1670.2
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2793.2
$ws.Range("I89").Value = 2793.2
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 13966
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -8350
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 31780
$ws.Range("J133").Value = 31780
$ws.Range("L133").Value = 31780
$ws.Range("N133").Value = -41900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1047.0785
$ws.Range("I58").Value = 1072.7805
$ws.Range("J58").Value = 941.7
$ws.Range("K58").Value = 1072.7805
$ws.Range("L58").Value = 941.7
$ws.Range("M58").Value = -869.7805000000001
$ws.Range("N58").Value = -1347.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 200000
$ws.Range("J80").Value = 200000
$ws.Range("L80").Value = 200000
$ws.Range("N80").Value = -202246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 200000
$ws.Range("J83").Value = 200000
$ws.Range("L83").Value = 600000
$ws.Range("N83").Value = -611232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1066.3462
$ws.Range("J107").Value = 326.63635
$ws.Range("L107").Value = 326.63635
$ws.Range("N107").Value = -4166.63635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1047.0785
$ws.Range("I136").Value = 1072.7805
$ws.Range("J136").Value = 941.7
$ws.Range("K136").Value = 3218.3415
$ws.Range("L136").Value = 2825.1
$ws.Range("M136").Value = -668.3415000000005
$ws.Range("N136").Value = -7925.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1178.1
$ws.Range("I132").Value = 704.4
$ws.Range("J132").Value = 1651.8
$ws.Range("K132").Value = 6339.599999999999
$ws.Range("L132").Value = 14866.2
$ws.Range("M132").Value = -3809.599999999999
$ws.Range("N132").Value = -19926.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 29118
$ws.Range("J15").Value = 29118
$ws.Range("L15").Value = 29118
$ws.Range("N15").Value = -29694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3381
$ws.Range("I80").Value = 2976.25
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 2976.25
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -1978.25
$ws.Range("N80").Value = -6996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 29118
$ws.Range("J81").Value = 29118
$ws.Range("L81").Value = 29118
$ws.Range("N81").Value = -31114

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3381
$ws.Range("I83").Value = 2976.25
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 14881.25
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -9889.25
$ws.Range("N83").Value = -34984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 29118
$ws.Range("J84").Value = 29118
$ws.Range("L84").Value = 87354
$ws.Range("N84").Value = -97338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 525
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 600
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 600
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -976

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3755.0151
$ws.Range("I132").Value = 4073.681
$ws.Range("J132").Value = 2966.7368
$ws.Range("K132").Value = 12221.043
$ws.Range("L132").Value = 8900.2104
$ws.Range("M132").Value = -9691.043
$ws.Range("N132").Value = -13960.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2352.9673
$ws.Range("I136").Value = 1795.1538
$ws.Range("J136").Value = 3341.818
$ws.Range("K136").Value = 5385.4614
$ws.Range("L136").Value = 10025.454
$ws.Range("M136").Value = -2835.4614
$ws.Range("N136").Value = -15125.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 12800
$ws.Range("J86").Value = 12800
$ws.Range("L86").Value = 12800
$ws.Range("N86").Value = -15046

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 12800
$ws.Range("J89").Value = 12800
$ws.Range("L89").Value = 64000
$ws.Range("N89").Value = -75232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4261.279
$ws.Range("I132").Value = 4513.7036
$ws.Range("J132").Value = 2314
$ws.Range("K132").Value = 13541.1108
$ws.Range("L132").Value = 6942
$ws.Range("M132").Value = -11011.1108
$ws.Range("N132").Value = -12002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5340.7407
$ws.Range("I136").Value = 5632.9165
$ws.Range("J136").Value = 3003.3333
$ws.Range("K136").Value = 16898.7495
$ws.Range("L136").Value = 9009.999899999999
$ws.Range("M136").Value = -14348.7495
$ws.Range("N136").Value = -14109.9999
